# NEL_sitrep_hospitalisation.xlsx update
# - Add daily hospitalisation figures for 2020-04-26 .. 2020-04-29 (rows 53-56)
# - Update the saved sheet view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data rows -------------------------------------------------------
# Date (col A), Barts admissions/discharges (B/C), Homerton (E/F), BHRUT (H/I)
# Columns D, G, J, K, L, M are all formulas mirroring the pattern used by
# the rows immediately above them.
$rows = @(
    @{ Row = 53; Date = 43947; B = 128; C = 47;  E = 315; F = 122; H = 49; I = 12 },
    @{ Row = 54; Date = 43948; B = 172; C = 50;  E = 312; F = 123; H = 48; I = 12 },
    @{ Row = 55; Date = 43949; B = 89;  C = 43;  E = 291; F = 113; H = 44; I = 12 },
    @{ Row = 56; Date = 43950; B = 133; C = 39;  E = 277; F = 114; H = 43; I = 12 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.Date
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Formula = "=B$n-C$n"
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Formula = "=E$n-F$n"
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Formula = "=H$n-I$n"
    $ws.Range("K$n").Formula = "=B$n+E$n+H$n"
    $ws.Range("L$n").Formula = "=C$n+F$n+I$n"
    $ws.Range("M$n").Formula = "=K$n-L$n"
}

# Rows 54-56 keep the pale-yellow highlight on column I that the sheet
# already uses for this block of manually-entered BHRUT_ITU figures.
$ws.Range("I54").Interior.ColorIndex = 6
$ws.Range("I55").Interior.ColorIndex = 6
$ws.Range("I56").Interior.ColorIndex = 6

# --- sheet view: scroll position + active selection -----------------------
$window = $excel.ActiveWindow
$window.ScrollRow = 31
$ws.Range("G48").Select()
